$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats
$xlPasteFormats = -4122

# Copy the existing style from A2 (used for column-A "label" cells) onto the
# new rows before we touch their values, so they pick up style index 1
# (same thin border / bold / centered-top alignment) instead of minting a
# new duplicate style.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A4").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A5").PasteSpecial($xlPasteFormats) | Out-Null

# Row 2: B2 value changes from 243 to 182 (A2 stays 0)
$ws.Range("B2").Value = 182

# Row 3 (new): A3=3, B3=109
$ws.Range("A3").Value = 3
$ws.Range("B3").Value = 109

# Row 4 (former row 3 shifted down): A4=1, B4=106
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = 106

# Row 5 (new): A5=2, B5=72
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = 72
